# Update symbol list (crypto prices) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values stored as text (inlineStr) in the
# source workbook. Prefixing with a leading apostrophe forces Excel to keep
# them as text instead of silently converting them to numbers.
function Set-TextValue($range, $text) {
    $ws.Range($range).Value = "'" + $text
}

Set-TextValue "D2"  "281.86"
Set-TextValue "D3"  "20.64"
Set-TextValue "D4"  "6.257"
Set-TextValue "D5"  "0.06146"
Set-TextValue "D7"  "6.564"
Set-TextValue "D8"  "1.504"
Set-TextValue "D9"  "0.8173"
Set-TextValue "D10" "0.01381"
Set-TextValue "D11" "0.1638"
Set-TextValue "D12" "0.08387"
Set-TextValue "D13" "0.03527"
Set-TextValue "D14" "0.03180"
Set-TextValue "D15" "0.09136"
Set-TextValue "D16" "3.703"
Set-TextValue "D17" "0.001640"
Set-TextValue "D18" "0.04703"
Set-TextValue "D19" "0.006542"
Set-TextValue "D20" "0.006158"
Set-TextValue "D21" "0.001067"
Set-TextValue "D22" "0.0001502"
Set-TextValue "D23" "3.770"
Set-TextValue "D40" "0.04681"
Set-TextValue "D41" "0.007161"

# Rows 42 and 43 swap their coin identity (BKEXToken <-> CEJI) and get
# refreshed price/volume-id data.
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.004505"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1098"
$ws.Range("E43").Value = "42BKEXTokenBKK"

Set-TextValue "D44" "0.01104"
Set-TextValue "D45" "0.00006526"

Write-Host "Applied cryptos.xlsx price refresh"
